$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 13
$ws.Range("G8").Value = 55.17
$ws.Range("H8").Value = 44.83
$ws.Range("J8").Value = 12
$ws.Range("K8").Value = 41.38
$ws.Range("E12").Value = 28
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 82.34999999999999
$ws.Range("H12").Value = 17.65
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 14.71
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 83.33
$ws.Range("H19").Value = 16.67
$ws.Range("I19").Value = 6.6
$ws.Range("J19").Value = 6
$ws.Range("K19").Value = 16.67
$ws.Range("E20").Value = 26
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 83.87
$ws.Range("H20").Value = 16.13
$ws.Range("I20").Value = 7.5
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 16.13
$ws.Range("E21").Value = 17
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 80.95
$ws.Range("H21").Value = 19.05
$ws.Range("I21").Value = 7.5
$ws.Range("J21").Value = 4
$ws.Range("K21").Value = 19.05
$ws.Range("E22").Value = 33
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 82.5
$ws.Range("H22").Value = 17.5
$ws.Range("I22").Value = 7.8
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 17.5
$ws.Range("E23").Value = 19
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 82.61
$ws.Range("H23").Value = 17.39
$ws.Range("I23").Value = 6.9
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 17.39

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 18
$ws.Range("G8").Value = 37.93
$ws.Range("H8").Value = 62.07
$ws.Range("I8").Value = 6.4
$ws.Range("J8").Value = 17
$ws.Range("K8").Value = 58.62
$ws.Range("E9").Value = 25
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 75.76000000000001
$ws.Range("H9").Value = 24.24
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 8
$ws.Range("K9").Value = 24.24
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = 19
$ws.Range("G10").Value = 50
$ws.Range("H10").Value = 50
$ws.Range("I10").Value = 8.5
$ws.Range("J10").Value = 19
$ws.Range("K10").Value = 50
$ws.Range("E11").Value = 25
$ws.Range("F11").Value = 14
$ws.Range("G11").Value = 64.09999999999999
$ws.Range("H11").Value = 35.9
$ws.Range("I11").Value = 8.5
$ws.Range("J11").Value = 14
$ws.Range("K11").Value = 35.9
$ws.Range("E12").Value = 24
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 70.59
$ws.Range("H12").Value = 29.41
$ws.Range("J12").Value = 10
$ws.Range("K12").Value = 29.41
$ws.Range("E13").Value = 29
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 74.36
$ws.Range("H13").Value = 25.64
$ws.Range("I13").Value = 7.2
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 90.23999999999999
$ws.Range("H14").Value = 9.76
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 23
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 92
$ws.Range("H15").Value = 8
$ws.Range("I15").Value = 7.3
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 33
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 84.62
$ws.Range("H16").Value = 15.38
$ws.Range("I16").Value = 7.2
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("E17").Value = 27
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 77.14
$ws.Range("H17").Value = 22.86
$ws.Range("I17").Value = 7.3
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("E18").Value = 27
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 75
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 7.4
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 83.33
$ws.Range("H19").Value = 16.67
$ws.Range("I19").Value = 6.9
$ws.Range("J19").Value = 6
$ws.Range("K19").Value = 16.67
$ws.Range("E20").Value = 26
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 83.87
$ws.Range("H20").Value = 16.13
$ws.Range("I20").Value = 7.5
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 16.13
$ws.Range("E21").Value = 17
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 80.95
$ws.Range("H21").Value = 19.05
$ws.Range("I21").Value = 7.8
$ws.Range("J21").Value = 4
$ws.Range("K21").Value = 19.05
$ws.Range("E22").Value = 33
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 82.5
$ws.Range("H22").Value = 17.5
$ws.Range("I22").Value = 7.5
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 17.5
$ws.Range("E23").Value = 19
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 82.61
$ws.Range("H23").Value = 17.39
$ws.Range("I23").Value = 6.8
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 17.39

$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 13
$ws.Range("G8").Value = 55.17
$ws.Range("H8").Value = 44.83
$ws.Range("I8").Value = 6.6
$ws.Range("J8").Value = 12
$ws.Range("K8").Value = 41.38
$ws.Range("E12").Value = 29
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 85.29000000000001
$ws.Range("H12").Value = 14.71
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 14.71
$ws.Range("I13").Value = 7.3
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 90.23999999999999
$ws.Range("H14").Value = 9.76
$ws.Range("I14").Value = 8.199999999999999
$ws.Range("E16").Value = 33
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 84.62
$ws.Range("H16").Value = 15.38
$ws.Range("I16").Value = 7.5
$ws.Range("E17").Value = 27
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 77.14
$ws.Range("H17").Value = 22.86
$ws.Range("I17").Value = 7.5
$ws.Range("E18").Value = 27
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 75
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 7.5
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 6
$ws.Range("G19").Value = 83.33
$ws.Range("H19").Value = 16.67
$ws.Range("I19").Value = 6.9
$ws.Range("J19").Value = 6
$ws.Range("K19").Value = 16.67
$ws.Range("E20").Value = 26
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 83.87
$ws.Range("H20").Value = 16.13
$ws.Range("I20").Value = 7.7
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 16.13
$ws.Range("E21").Value = 17
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 80.95
$ws.Range("H21").Value = 19.05
$ws.Range("I21").Value = 7.9
$ws.Range("J21").Value = 4
$ws.Range("K21").Value = 19.05
$ws.Range("E22").Value = 33
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 82.5
$ws.Range("H22").Value = 17.5
$ws.Range("I22").Value = 7.8
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 17.5
$ws.Range("E23").Value = 19
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 82.61
$ws.Range("H23").Value = 17.39
$ws.Range("I23").Value = 7.1
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 17.39

